$wb = $excel.ActiveWorkbook

# Rename the worksheet-connection defined names so each gets a trailing "1"
# (supports the new ability to reference a *collection* of ranges in the
# FROM_EXCEL transformer, introduced by this change).
foreach ($n in $wb.Names) {
    $n.Name = $n.Name + "1"
}

# Add an extra data row to the RepeatingData sheet used by the new tests.
$ws3 = $wb.Worksheets.Item("RepeatingData")
$ws3.Range("A9").Value = "name1"
$ws3.Range("B9").Value = "status1"

# Switch the active sheet/selection to RepeatingData (matches the workbook
# being left open on that sheet after the edits were made).
$ws3.Activate()
$ws3.Range("A5").Select()
